$wb = $excel.ActiveWorkbook

# --- Task 1: it is no longer the active/selected tab ---
$wsTask1 = $wb.Worksheets.Item("Task 1")

# --- Task 2: selection changes from P30 to B5:B13 ---
$wsTask2 = $wb.Worksheets.Item("Task 2")
$wsTask2.Activate()
$wsTask2.Range("B5:B13").Select()

# --- Task 3: becomes the active/selected tab, header + data updated ---
$wsTask3 = $wb.Worksheets.Item("Task 3")
$wsTask3.Activate()

# Update header B4 from "fn in Hz" to "f1 in Hz" (with "1" as subscript)
$headerCell = $wsTask3.Range("B4")
$headerCell.Value = "f1 in Hz"
$headerChars = $headerCell.Characters(2, 1)
$headerChars.Font.Subscript = $true
$restChars = $headerCell.Characters(3, 6)
$restChars.Font.Subscript = $false

# Update the fundamental-frequency measurements in column B (rows 5-14)
$wsTask3.Range("B5").Value = 161
$wsTask3.Range("B6").Value = 325
$wsTask3.Range("B7").Value = 489
$wsTask3.Range("B8").Value = 654
$wsTask3.Range("B9").Value = 818
$wsTask3.Range("B10").Value = 984
$wsTask3.Range("B11").Value = 1150
$wsTask3.Range("B12").Value = 1314
$wsTask3.Range("B13").Value = 1479
$wsTask3.Range("B14").Value = 1800

# Final selection on Task 3 is the single cell B5
$wsTask3.Range("B5").Select()
